$d = $word.ActiveDocument

# The document starts with an empty paragraph, then a paragraph holding the
# lab image, then a trailing empty paragraph. Turn the first (empty)
# paragraph into the two centered title lines.
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = "ACS 54500 - Lab 1 `rBy-Vikramadithya Ivaturi"

# Center both new title paragraphs.
$d.Paragraphs.Item(1).Alignment = 1
$d.Paragraphs.Item(2).Alignment = 1

# Mark the picture's run as NoProof (adds <w:noProof/> to its rPr), as in
# the original author's edit.
$picturePara = $d.Paragraphs.Item(3)
$picturePara.Range.NoProofing = 1
